# Update "想去人数" (F column) values on the "展览" and "全部类型" worksheets
# to reflect freshly generated data (gh-pages output at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (sheet1) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 237
$ws1.Range("F3").Value  = 437
$ws1.Range("F4").Value  = 13131
$ws1.Range("F6").Value  = 223
$ws1.Range("F13").Value = 70
$ws1.Range("F17").Value = 422
$ws1.Range("F18").Value = 5563
$ws1.Range("F22").Value = 14
$ws1.Range("F24").Value = 132
$ws1.Range("F25").Value = 155

# --- Sheet "全部类型" (sheet4) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 237
$ws4.Range("F3").Value  = 437
$ws4.Range("F4").Value  = 13132
$ws4.Range("F6").Value  = 223
$ws4.Range("F13").Value = 70
$ws4.Range("F17").Value = 422
$ws4.Range("F18").Value = 5563
$ws4.Range("F22").Value = 14
$ws4.Range("F24").Value = 132
$ws4.Range("F25").Value = 155
